# Insert a new "Net Income" line item between "Common Stock" (row 10) and
# "Retained Earnings" (row 11) on the balance sheet. This shifts the
# existing Retained Earnings / Total Equity / Total Liabilities & Equity
# rows down by one, and updates their figures accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11:13 down to 12:14 to make room for the new row.
$ws.Rows.Item(11).Insert()

# The figures are stored as literal text (not numbers) in the source
# workbook. A leading apostrophe forces Excel to keep the value as text
# (quotePrefix) instead of auto-converting "1,525" into a number.

# New row 11: Net Income
$ws.Range("A11").Value = "  Net Income"
$ws.Range("B11").Value = "net_income"
$ws.Range("C11").Value = "'1,525"
$ws.Range("D11").Value = "'1,850"
$ws.Range("E11").Value = "'2,116"
$ws.Range("F11").Value = "'2,429"
$ws.Range("G11").Value = "'2,799"

# Row 12: Retained Earnings (updated figures)
$ws.Range("A12").Value = "  Retained Earnings"
$ws.Range("B12").Value = "retained_earnings"
$ws.Range("C12").Value = "'1,635"
$ws.Range("D12").Value = "'2,905"
$ws.Range("E12").Value = "'4,316"
$ws.Range("F12").Value = "'5,984"
$ws.Range("G12").Value = "'7,982"

# Row 13: Total Equity (updated figures)
$ws.Range("A13").Value = "  Total Equity"
$ws.Range("B13").Value = "total_equity"
$ws.Range("C13").Value = "'1,735"
$ws.Range("D13").Value = "'3,005"
$ws.Range("E13").Value = "'4,416"
$ws.Range("F13").Value = "'6,084"
$ws.Range("G13").Value = "'8,082"

# Row 14: Total Liabilities & Equity (updated figures)
$ws.Range("A14").Value = "  Total Liabilities & Equity"
$ws.Range("B14").Value = "total_liabs_equity"
$ws.Range("C14").Value = "'1,965"
$ws.Range("D14").Value = "'3,235"
$ws.Range("E14").Value = "'4,656"
$ws.Range("F14").Value = "'6,333"
$ws.Range("G14").Value = "'8,342"
